$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.076.45"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "2.365.56"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.62"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("E6").Value = "  -3.84%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "2.372.36"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.76"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.323"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "2.786.63"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").Value = "56.027.95"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.41"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "2.282.83"
$ws.Range("E18").Value = "  -6.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.99"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.23"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.28"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.85"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.368"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.147"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("E28").Value = "  -5.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.21"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "0.0₃0712"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  -7.13%  "
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.59"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.07"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -6.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.50"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.67"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -6.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.562"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "239.83"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -6.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0479"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.03"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.950"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.62%  "
